# Update COVID-19 country stats table ("Pais" sheet) to the later snapshot.
#
# Column A (country) cells reference shared strings by index, and several
# rows' A-cells were re-pointed to a *different* existing country name
# (i.e. the country list around those rows got reshuffled) while columns
# B-H carry refreshed case/death counts. For every reshuffled group we
# first stash each cell behind a unique temp placeholder before writing
# final values — writing straight to final values risks two cells
# momentarily holding onto the same shared-string text and Excel
# collapsing them onto one shared-string slot (so both cells would end up
# showing the same text instead of swapping).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Junio de 2020 a las 14:43"

# --- Reorder country names (swap/rotate via temp placeholders to avoid shared-string collisions) ---
$ws.Range("A76").Value = "__TMP_76__"
$ws.Range("A77").Value = "__TMP_77__"
$ws.Range("A76").Value = "Uzbekistan"
$ws.Range("A77").Value = "Tayikistan"

$ws.Range("A90").Value = "__TMP_90__"
$ws.Range("A91").Value = "__TMP_91__"
$ws.Range("A92").Value = "__TMP_92__"
$ws.Range("A90").Value = "Etiopia"
$ws.Range("A91").Value = "Tailandia"
$ws.Range("A92").Value = "Grecia"

$ws.Range("A105").Value = "__TMP_105__"
$ws.Range("A106").Value = "__TMP_106__"
$ws.Range("A105").Value = "Mali"
$ws.Range("A106").Value = "Lituania"

$ws.Range("A143").Value = "__TMP_143__"
$ws.Range("A144").Value = "__TMP_144__"
$ws.Range("A145").Value = "__TMP_145__"
$ws.Range("A143").Value = "Mozambique"
$ws.Range("A144").Value = "Togo"
$ws.Range("A145").Value = "Ruanda"

$ws.Range("A206").Value = "__TMP_206__"
$ws.Range("A207").Value = "__TMP_207__"
$ws.Range("A206").Value = "Groenlandia"
$ws.Range("A207").Value = "Islas Malvinas"

$ws.Range("A208").Value = "__TMP_208__"
$ws.Range("A209").Value = "__TMP_209__"
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("A209").Value = "Islas Turcas y Caicos"

$ws.Range("A210").Value = "__TMP_210__"
$ws.Range("A211").Value = "__TMP_211__"
$ws.Range("A210").Value = "Montserrat"
$ws.Range("A211").Value = "Seychelles"

$ws.Range("A213").Value = "__TMP_213__"
$ws.Range("A214").Value = "__TMP_214__"
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("A214").Value = "Papua Nueva Guinea"

# --- Update numeric statistics ---
$ws.Range("B4").Value = 2117333
$ws.Range("C4").Value = 411
$ws.Range("D4").Value = 842006
$ws.Range("E4").Value = 1158487
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 116840

$ws.Range("B5").Value = 831064
$ws.Range("C5").Value = 1162
$ws.Range("E5").Value = 361502
$ws.Range("G5").Value = 51
$ws.Range("H5").Value = 41952

$ws.Range("B7").Value = 310367
$ws.Range("C7").Value = 764
$ws.Range("D7").Value = 155290
$ws.Range("E7").Value = 146182

$ws.Range("B12").Value = 187256
$ws.Range("C12").Value = 5
$ws.Range("E12").Value = 6493

$ws.Range("B27").Value = 50931
$ws.Range("C27").Value = 180
$ws.Range("G27").Value = 20
$ws.Range("H27").Value = 4874

$ws.Range("B28").Value = 48640
$ws.Range("C28").Value = 179
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 6057

$ws.Range("B31").Value = 41990
$ws.Range("C31").Value = 491
$ws.Range("D31").Value = 26761
$ws.Range("E31").Value = 14941
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 288

$ws.Range("B36").Value = 35466
$ws.Range("C36").Value = 514
$ws.Range("D36").Value = 25882
$ws.Range("E36").Value = 9295
$ws.Range("G36").Value = 4
$ws.Range("H36").Value = 289

$ws.Range("B58").Value = 12139
$ws.Range("C58").Value = 40
$ws.Range("D58").Value = 11035
$ws.Range("E58").Value = 507
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 597

$ws.Range("B76").Value = 4937
$ws.Range("C76").Value = 68
$ws.Range("D76").Value = 3837
$ws.Range("E76").Value = 1081
$ws.Range("H76").Value = 19

$ws.Range("B77").Value = 4902
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 3158
$ws.Range("E77").Value = 1695
$ws.Range("H77").Value = 49

$ws.Range("B90").Value = 3166
$ws.Range("C90").Value = 251
$ws.Range("D90").Value = 495
$ws.Range("E90").Value = 2616
$ws.Range("G90").Value = 8
$ws.Range("H90").Value = 55

$ws.Range("B91").Value = 3134
$ws.Range("C91").Value = 5
$ws.Range("D91").Value = 2987
$ws.Range("E91").Value = 89
$ws.Range("H91").Value = 58

$ws.Range("B92").Value = 3108
$ws.Range("D92").Value = 1374
$ws.Range("E92").Value = 1551
$ws.Range("H92").Value = 183

$ws.Range("B97").Value = 2251
$ws.Range("C97").Value = 2
$ws.Range("D97").Value = 2134
$ws.Range("E97").Value = 10

$ws.Range("B103").Value = 1882
$ws.Range("C103").Value = 2
$ws.Range("E103").Value = 619

$ws.Range("B105").Value = 1776
$ws.Range("C105").Value = 24
$ws.Range("D105").Value = 1058
$ws.Range("E105").Value = 614
$ws.Range("G105").Value = 3
$ws.Range("H105").Value = 104

$ws.Range("B106").Value = 1763
$ws.Range("C106").Value = 7
$ws.Range("D106").Value = 1416
$ws.Range("E106").Value = 272
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 75

$ws.Range("B116").Value = 1442
$ws.Range("C116").Value = 20
$ws.Range("D116").Value = 868
$ws.Range("E116").Value = 542
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 32

$ws.Range("B121").Value = 1110
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 1061

$ws.Range("B143").Value = 553
$ws.Range("C143").Value = 44
$ws.Range("D143").Value = 151
$ws.Range("E143").Value = 400
$ws.Range("H143").Value = 2

$ws.Range("B144").Value = 525
$ws.Range("D144").Value = 279
$ws.Range("E144").Value = 233
$ws.Range("H144").Value = 13

$ws.Range("B145").Value = 510
$ws.Range("D145").Value = 321
$ws.Range("E145").Value = 187

$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0

